$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write values in the same order the original author must have entered them so
# that new shared-string entries land at the expected indices (65..69):
#   65 tipo
#   66 FONDOS PUBL. ADMINIST. PENSIONES
#   67 FONDOS PRIV. ADMINIST. PENSIONES
#   68 EMP. PRIV. PROMOTORAS DE SALUD
#   69 PARAFISCALES
$ws.Range("I1").Value = "tipo"
$ws.Range("I15").Value = "FONDOS PUBL. ADMINIST. PENSIONES"
$ws.Range("I16").Value = "FONDOS PRIV. ADMINIST. PENSIONES"
$ws.Range("I2").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I20").Value = "PARAFISCALES"

# Match the header cell's style (bold) by copying the format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Fill in the rest of the "tipo" column.
$ws.Range("I3").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I4").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I5").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I6").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I7").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I8").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I9").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I10").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I11").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I12").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I13").Value = "EMP. PRIV. PROMOTORAS DE SALUD"
$ws.Range("I14").Value = "EMP. PRIV. PROMOTORAS DE SALUD"

$ws.Range("I17").Value = "FONDOS PRIV. ADMINIST. PENSIONES"
$ws.Range("I18").Value = "FONDOS PRIV. ADMINIST. PENSIONES"

$ws.Range("I19").Value = "RIESGOS PROFESIONALES"

$ws.Range("I21").Value = "PARAFISCALES"
$ws.Range("I22").Value = "PARAFISCALES"
$ws.Range("I23").Value = "PARAFISCALES"
$ws.Range("I24").Value = "PARAFISCALES"

# Widen the new column, mirroring the author's manual resize.
$ws.Range("I1").ColumnWidth = 43

# Leave the selection where the author left it when they saved.
$ws.Range("G28").Select()
